$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 (was "Top Wear"/"Cloths" -> now "New Fashions"/"Fashion") ---
$ws.Range("A5").Value = "67139ea4c9907d0ec0ce54c7"
$ws.Range("B5").Value = "New Fashions"
$ws.Range("C5").Value = "Fashion"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 45584.72748332176
$ws.Range("F5").Value = 45584.72748332176

# --- New row 6 ---
$ws.Range("A6").Value = "672a18601c261e0e7b03e778"
$ws.Range("B6").Value = "Healty"
$ws.Range("C6").Value = "Fruits"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 45601.77558584491
$ws.Range("F6").Value = 45601.77558584491

# --- New row 7 ---
$ws.Range("A7").Value = "672a43c0004d25297f6a4908"
$ws.Range("B7").Value = "banana"
$ws.Range("C7").Value = "Fruits"
$ws.Range("D7").Value = "medium"
$ws.Range("E7").Value = 45601.904101238426
$ws.Range("F7").Value = 45601.904101238426

# --- New row 8 ---
$ws.Range("A8").Value = "672a4435004d25297f6a4917"
$ws.Range("B8").Value = "Laptops"
$ws.Range("C8").Value = "Electronics"
$ws.Range("D8").Value = "low"
$ws.Range("E8").Value = 45601.90546107639
$ws.Range("F8").Value = 45601.911041493055

# --- New row 9 ---
$ws.Range("A9").Value = "672a44ed004d25297f6a493e"
$ws.Range("B9").Value = "Mobiles"
$ws.Range("C9").Value = "Electronics"
$ws.Range("D9").Value = "low"
$ws.Range("E9").Value = 45601.90758197917
$ws.Range("F9").Value = 45601.955005497686

# The date columns (E/F) carry a date number-format style (style index 1 in the
# original sheet, applied to E2:F5). Copy that formatting down onto the newly
# added rows so E6:F9 pick up the same style instead of a freshly-minted one.
$ws.Range("E2:F2").Copy() | Out-Null
$ws.Range("E6:F9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
